$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "AÇÂO" (action) texts for the remaining risk rows (H7:H10)
$ws.Range("H7").Value = "         Manter o grupo todo atualizado com os requisitos pendentes e com as datas das entregas."
$ws.Range("H8").Value = "Fazer reunião emergencial, discutir metas pendentes e fazer a redivisão das tarefas."
$ws.Range("H9").Value = "Explicar tudo o que está a ser feito, manter o grupo alinhado em relação ao projeto e procurar ajuda dos membros "
$ws.Range("H10").Value = "Dividir as atividades de acordo com as capacidades individuais dos integrantes do grupo"

# Match the formatting of the already-filled action cells (bold, centered both ways)
# by copying the format from H5 (which already has the target style) onto H4 and H7:H10.
$ws.Range("H5").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H7:H10").PasteSpecial(-4122)

# Restore the active cell selection to match the saved workbook.
$ws.Range("H14").Select()
